# Update countries & provincias Spain
# - Update "last updated" timestamp (09:35 -> 10:05)
# - Refresh several countries' COVID figures (Rusia, Estonia, Lituania, Eslovaquia)
# - Reorder Letonia/Maldivas (Letonia now comes before Maldivas) and refresh their figures
# - Refresh Sri Lanka figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 10:05"

# Row 6: Rusia
$ws.Range("B6").Value = 252245
$ws.Range("C6").Value = 9974
$ws.Range("D6").Value = 53530
$ws.Range("E6").Value = 196410
$ws.Range("G6").Value = 93
$ws.Range("H6").Value = 2305

# Row 86: Estonia
$ws.Range("B86").Value = 1758
$ws.Range("C86").Value = 7
$ws.Range("D86").Value = 909
$ws.Range("E86").Value = 787
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 62

# Row 88: Lituania
$ws.Range("B88").Value = 1511
$ws.Range("C88").Value = 6
$ws.Range("D88").Value = 934
$ws.Range("E88").Value = 523

# Row 90: Eslovaquia
$ws.Range("B90").Value = 1477
$ws.Range("C90").Value = 8
$ws.Range("D90").Value = 1112
$ws.Range("E90").Value = 338
$ws.Range("F90").Value = 5

# Rows 102-103: swap Maldivas/Letonia order and refresh figures.
# Row 102 becomes Letonia (new figures), row 103 becomes Maldivas (old row 102 figures)
$ws.Range("A102").Value = "Letonia"
$ws.Range("B102").Value = 962
$ws.Range("C102").Value = 11
$ws.Range("D102").Value = 627
$ws.Range("E102").Value = 316
$ws.Range("H102").Value = 19

$ws.Range("A103").Value = "Maldivas"
$ws.Range("B103").Value = 955
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 40
$ws.Range("E103").Value = 911
$ws.Range("H103").Value = 4

# Row 104: Sri Lanka
$ws.Range("D104").Value = 445
$ws.Range("E104").Value = 461
